# Data Drive Script is added
# Adds a new "ValidLogin" worksheet after the existing "TC1" sheet,
# populates it with username/password test data, and makes it the
# active/selected sheet.

$wb = $excel.ActiveWorkbook

# Existing sheet (TC1) - the new sheet should be inserted right after it.
$tc1 = $wb.Worksheets.Item(1)

# Insert the new worksheet right after TC1.
$newSheet = $wb.Worksheets.Add($null, $tc1)
$newSheet.Name = "ValidLogin"

# Populate the data-drive table.
$newSheet.Range("A1").Value = "username"
$newSheet.Range("B1").Value = "password"
$newSheet.Range("A2").Value = "admin"
$newSheet.Range("B2").Value = "pointofsale"

# Match the authored selection / zoom on the new (now active) sheet.
$newSheet.Range("B3").Select()
$excel.ActiveWindow.Zoom = 160
